# Apply the "jan 15 2024 timepix command document" revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- K column (reply instruction hex) updates, rows 4-9 ---
$ws.Range("K4").Value = "0x1b"
$ws.Range("K5").Value = "0x09"
$ws.Range("K6").Value = "0x01"
$ws.Range("K7").Value = "0x08"
$ws.Range("K8").Value = "0x08"
$ws.Range("K9").Value = "0x04"

# --- A7 / A8 label font weight: bold -> normal (matches rows 16-19) ---
$ws.Range("A7").Font.Bold = $false
$ws.Range("A8").Font.Bold = $false

# --- read_rates (row 9) bitstring value changes: 101000 -> 1 ---
# (E9 / AA9 formulas recompute automatically: 0xA8 -> 0x81)
$ws.Range("D9").Value = 1

# --- erase_storage (row 13) bitstring value changes: 10001 -> 100010 ---
# (E13 / AA13 formulas recompute automatically: 0x11 -> 0x22)
$ws.Range("D13").Value = 100010

# --- applicability flags (rows 11-19): clear Formatter/Housekeeping and the
#     CdTe1-CdTeDE / CMOS1 / CMOS2 columns (L,M,P,Q,R,S,T,U,V) to 0, leaving
#     GSE/EVTM (N,O) and Timepix (W) untouched ---
$rows = 11..19
foreach ($r in $rows) {
    $ws.Range("L" + $r).Value = 0
    $ws.Range("M" + $r).Value = 0
    $ws.Range("P" + $r).Value = 0
    $ws.Range("Q" + $r).Value = 0
    $ws.Range("R" + $r).Value = 0
    $ws.Range("S" + $r).Value = 0
    $ws.Range("T" + $r).Value = 0
    $ws.Range("U" + $r).Value = 0
    $ws.Range("V" + $r).Value = 0
}

# --- sheet view: restore the selection the author left behind ---
[void]$ws.Range("V24").Select()

Write-Host "edit complete"
